# Apply the "Iteration 5" round update to the "Ghi cong" (scoreboard) sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ghi cong")

# Header for the evaluation round that used to be labelled "ĐG lần 6" (no
# date) now gets its date stamp, matching the other round headers.
$ws.Range("K5").Value = "ĐG lần 6" + [char]10 + "(24/05)"

# New scores recorded for the "ĐG lần 5 (20/05)" round (column J) for the
# five students in rows 11-15.
$ws.Cells.Item(11, 10).Value = 1
$ws.Cells.Item(12, 10).Value = 1
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(14, 10).Value = 1
$ws.Cells.Item(15, 10).Value = 0

# Move the active selection like the author's saved view.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("K6").Select()
